# Update NATMI Wnt5a-Fzd5 LR-pair sheet with recalculated TPM-based values.
# The sending-cluster set grew from {FAPs} to {ECs, FAPs}, so every
# (sending cluster x target cluster) combination is rewritten and six new
# rows (ECs as sender) are appended below the existing six (FAPs as sender).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt5a"
$ws.Range("C2").Value = "Fzd5"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.125615
$ws.Range("H2").Value = 0.25123
$ws.Range("I2").Value = 0.02647478672532295
$ws.Range("J2").Value = 0.01780700335556722
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 2.8377025
$ws.Range("N2").Value = 5.675405
$ws.Range("O2").Value = 0.08520139853031897
$ws.Range("P2").Value = 0.05979321158534227
$ws.Range("Q2").Value = 0.3564579995375
$ws.Range("R2").Value = 1.42583199815
$ws.Range("S2").Value = 0.002255688854789439
$ws.Range("T2").Value = 0.00106473791934033

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt5a"
$ws.Range("C3").Value = "Fzd5"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.125615
$ws.Range("H3").Value = 0.25123
$ws.Range("I3").Value = 0.02647478672532295
$ws.Range("J3").Value = 0.01780700335556722
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.061974333333333
$ws.Range("N3").Value = 15.185923
$ws.Range("O3").Value = 0.1519846751111432
$ws.Range("P3").Value = 0.1599912441592654
$ws.Range("Q3").Value = 0.6358599058816666
$ws.Range("R3").Value = 3.81515943529
$ws.Range("S3").Value = 0.004023761859085017
$ws.Range("T3").Value = 0.002848964621605413

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt5a"
$ws.Range("C4").Value = "Fzd5"
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.125615
$ws.Range("H4").Value = 0.25123
$ws.Range("I4").Value = 0.02647478672532295
$ws.Range("J4").Value = 0.01780700335556722
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.372280666666666
$ws.Range("N4").Value = 13.116842
$ws.Range("O4").Value = 0.1312767732230829
$ws.Range("P4").Value = 0.1381924477702479
$ws.Range("Q4").Value = 0.5492240359433332
$ws.Range("R4").Value = 3.29534421566
$ws.Range("S4").Value = 0.003475524573069707
$ws.Range("T4").Value = 0.002460793381158852

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt5a"
$ws.Range("C5").Value = "Fzd5"
$ws.Range("D5").Value = "MuSCs"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.125615
$ws.Range("H5").Value = 0.25123
$ws.Range("I5").Value = 0.02647478672532295
$ws.Range("J5").Value = 0.01780700335556722
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.1625465
$ws.Range("N5").Value = 4.325093
$ws.Range("O5").Value = 0.06492998691259792
$ws.Range("P5").Value = 0.04556700374251401
$ws.Range("Q5").Value = 0.2716482785975
$ws.Range("R5").Value = 1.08659311439
$ws.Range("S5").Value = 0.001719007555589041
$ws.Range("T5").Value = 0.000811411788546091

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt5a"
$ws.Range("C6").Value = "Fzd5"
$ws.Range("D6").Value = "Neutrophils"
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.125615
$ws.Range("H6").Value = 0.25123
$ws.Range("I6").Value = 0.02647478672532295
$ws.Range("J6").Value = 0.01780700335556722
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 15.97657333333333
$ws.Range("N6").Value = 47.92972
$ws.Range("O6").Value = 0.479693129114909
$ws.Range("P6").Value = 0.5049634148023289
$ws.Range("Q6").Value = 2.006897259266667
$ws.Range("R6").Value = 12.0413835556
$ws.Range("S6").Value = 0.01269977328692002
$ws.Range("T6").Value = 0.008991885221823752

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Wnt5a"
$ws.Range("C7").Value = "Fzd5"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.125615
$ws.Range("H7").Value = 0.25123
$ws.Range("I7").Value = 0.02647478672532295
$ws.Range("J7").Value = 0.01780700335556722
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.894743333333333
$ws.Range("N7").Value = 8.68423
$ws.Range("O7").Value = 0.08691403710794816
$ws.Range("P7").Value = 0.0914926779403015
$ws.Range("Q7").Value = 0.3636231838166666
$ws.Range("R7").Value = 2.1817391029
$ws.Range("S7").Value = 0.002301030595869733
$ws.Range("T7").Value = 0.00162921042309278

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt5a"
$ws.Range("C8").Value = "Fzd5"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 4.619088000000001
$ws.Range("H8").Value = 13.857264
$ws.Range("I8").Value = 0.9735252132746771
$ws.Range("J8").Value = 0.9821929966444328
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 2.8377025
$ws.Range("N8").Value = 5.675405
$ws.Range("O8").Value = 0.08520139853031897
$ws.Range("P8").Value = 0.05979321158534227
$ws.Range("Q8").Value = 13.10759756532
$ws.Range("R8").Value = 78.64558539192
$ws.Range("S8").Value = 0.08294570967552953
$ws.Range("T8").Value = 0.05872847366600194

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt5a"
$ws.Range("C9").Value = "Fzd5"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 4.619088000000001
$ws.Range("H9").Value = 13.857264
$ws.Range("I9").Value = 0.9735252132746771
$ws.Range("J9").Value = 0.9821929966444328
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.061974333333333
$ws.Range("N9").Value = 15.185923
$ws.Range("O9").Value = 0.1519846751111432
$ws.Range("P9").Value = 0.1599912441592654
$ws.Range("Q9").Value = 23.381704899408
$ws.Range("R9").Value = 210.435344094672
$ws.Range("S9").Value = 0.1479609132520582
$ws.Range("T9").Value = 0.15714227953766

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt5a"
$ws.Range("C10").Value = "Fzd5"
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 4.619088000000001
$ws.Range("H10").Value = 13.857264
$ws.Range("I10").Value = 0.9735252132746771
$ws.Range("J10").Value = 0.9821929966444328
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 4.372280666666666
$ws.Range("N10").Value = 13.116842
$ws.Range("O10").Value = 0.1312767732230829
$ws.Range("P10").Value = 0.1381924477702479
$ws.Range("Q10").Value = 20.195949160032
$ws.Range("R10").Value = 181.763542440288
$ws.Range("S10").Value = 0.1278012486500132
$ws.Range("T10").Value = 0.1357316543890891

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt5a"
$ws.Range("C11").Value = "Fzd5"
$ws.Range("D11").Value = "MuSCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 4.619088000000001
$ws.Range("H11").Value = 13.857264
$ws.Range("I11").Value = 0.9735252132746771
$ws.Range("J11").Value = 0.9821929966444328
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.1625465
$ws.Range("N11").Value = 4.325093
$ws.Range("O11").Value = 0.06492998691259792
$ws.Range("P11").Value = 0.04556700374251401
$ws.Range("Q11").Value = 9.988992587592001
$ws.Range("R11").Value = 59.933955525552
$ws.Range("S11").Value = 0.06321097935700888
$ws.Range("T11").Value = 0.04475559195396792

# Row 12
$ws.Range("A12").Value = "FAPs"
$ws.Range("B12").Value = "Wnt5a"
$ws.Range("C12").Value = "Fzd5"
$ws.Range("D12").Value = "Neutrophils"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 4.619088000000001
$ws.Range("H12").Value = 13.857264
$ws.Range("I12").Value = 0.9735252132746771
$ws.Range("J12").Value = 0.9821929966444328
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 15.97657333333333
$ws.Range("N12").Value = 47.92972
$ws.Range("O12").Value = 0.479693129114909
$ws.Range("P12").Value = 0.5049634148023289
$ws.Range("Q12").Value = 73.79719816512001
$ws.Range("R12").Value = 664.1747834860801
$ws.Range("S12").Value = 0.466993355827989
$ws.Range("T12").Value = 0.4959715295805052

# Row 13
$ws.Range("A13").Value = "FAPs"
$ws.Range("B13").Value = "Wnt5a"
$ws.Range("C13").Value = "Fzd5"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 4.619088000000001
$ws.Range("H13").Value = 13.857264
$ws.Range("I13").Value = 0.9735252132746771
$ws.Range("J13").Value = 0.9821929966444328
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 2.894743333333333
$ws.Range("N13").Value = 8.68423
$ws.Range("O13").Value = 0.08691403710794816
$ws.Range("P13").Value = 0.0914926779403015
$ws.Range("Q13").Value = 13.37107419408
$ws.Range("R13").Value = 120.33966774672
$ws.Range("S13").Value = 0.08461300651207844
$ws.Range("T13").Value = 0.08986346751720872
